$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New row 16: quarterly date label "01-07-2021" plus the five index values.
# Force the date-looking label to be entered as text (matching the other
# "Serie" column entries) instead of letting Excel auto-convert it to a
# date serial number, then clear the temporary text format so the cell
# keeps the default (unstyled) appearance used by the rest of the column.
$dateCell = $ws.Range("A16")
$dateCell.NumberFormat = "@"
$dateCell.Value = "01-07-2021"
$dateCell.ClearFormats()

$ws.Range("B16").Value = 110.45
$ws.Range("C16").Value = 108.62
$ws.Range("D16").Value = 112.15
$ws.Range("E16").Value = 108.45
$ws.Range("F16").Value = 119.48
